$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.03426924709026265
$ws.Range("H2").Value = -18.94433775351375
$ws.Range("I2").Value = 13.26394738455261
$ws.Range("G3").Value = 0.0595546252267164
$ws.Range("H3").Value = 22.28788129378635
$ws.Range("G4").Value = -0.01635275640033013
$ws.Range("H4").Value = -958.9103726605158
$ws.Range("G5").Value = -0.01297376065785644
$ws.Range("H5").Value = -441.7429242367151
$ws.Range("G6").Value = 0.03504527399398531
$ws.Range("H6").Value = 1.083643191057361
$ws.Range("G7").Value = 0.06985041081686504
$ws.Range("H7").Value = 31.32505056207568
$ws.Range("G8").Value = 0.008706227010006257
$ws.Range("H8").Value = 146.2585998641752
$ws.Range("G9").Value = 0.0009685069108547311
$ws.Range("H9").Value = 104.498844222027
$ws.Range("G10").Value = -0.08467300880605795
$ws.Range("H10").Value = -16.46689066515937
$ws.Range("G11").Value = -0.06667102992806373
$ws.Range("H11").Value = 27.54018559148611
$ws.Range("G12").Value = -0.2552254454605632
$ws.Range("H12").Value = -4.406053015563982
$ws.Range("G13").Value = -0.3007508121401353
$ws.Range("H13").Value = -9.438261786552655
$ws.Range("G14").Value = -0.02421180379573063
$ws.Range("H14").Value = 34.73413635701889
$ws.Range("G15").Value = -0.03211262875972212
$ws.Range("H15").Value = 7.646323293681968
$ws.Range("G16").Value = 0.1101152374305845
$ws.Range("H16").Value = -12.12446714627186
$ws.Range("G17").Value = 0.1350710642820757
$ws.Range("H17").Value = -3.692983073542735
$ws.Range("G18").Value = 0.1363639852945159
$ws.Range("H18").Value = 9.324947708638275
$ws.Range("G19").Value = 0.109637546754897
$ws.Range("H19").Value = -17.69808273209284
$ws.Range("G20").Value = 0.0399523481137436
$ws.Range("H20").Value = 16.35664306429763
$ws.Range("G21").Value = 0.04701411861497632
$ws.Range("H21").Value = -18.9968471706373
$ws.Range("G22").Value = -0.09306538432034084
$ws.Range("H22").Value = -16.55320463522862
$ws.Range("G23").Value = -0.08525911434974869
$ws.Range("H23").Value = -36.31330889980579
$ws.Range("G24").Value = 0.1203945863661415
$ws.Range("H24").Value = 1.930394773516972
$ws.Range("G25").Value = 0.1340394161425693
$ws.Range("H25").Value = 6.23682162503695
$ws.Range("G26").Value = 0.03775439059890311
$ws.Range("H26").Value = -24.04171010468165
$ws.Range("G27").Value = 0.05830829231337434
$ws.Range("H27").Value = -32.72888843074746
$ws.Range("G28").Value = -0.08072515062216395
$ws.Range("H28").Value = -26.93823993925693
$ws.Range("G29").Value = -0.07599121292847734
$ws.Range("H29").Value = -6.770888596030828
$ws.Range("G30").Value = 0.05835898897489881
$ws.Range("H30").Value = -8.396319415411336
$ws.Range("G31").Value = 0.05350352906587293
$ws.Range("H31").Value = -11.68187242392493
$ws.Range("G32").Value = 0.08388967135687306
$ws.Range("H32").Value = -14.6282618287856
$ws.Range("G33").Value = 0.09785141731653284
$ws.Range("H33").Value = 18.91811140061141
$ws.Range("G34").Value = 0.002797348339994649
$ws.Range("H34").Value = -89.26376836130561
$ws.Range("G35").Value = -0.004577573919160895
$ws.Range("H35").Value = 59.16566052260834
$ws.Range("G36").Value = 0.008883952914180059
$ws.Range("H36").Value = 1550.657009402615
$ws.Range("G37").Value = -0.01178387172396554
$ws.Range("H37").Value = 6.136505967891603
$ws.Range("G38").Value = 0.1158172361505888
$ws.Range("H38").Value = 7.980782023990553
$ws.Range("G39").Value = 0.1010639943828278
$ws.Range("H39").Value = 17.97888424902274
$ws.Range("G40").Value = 0.01123280623663582
$ws.Range("H40").Value = 278.1762032314729
$ws.Range("G41").Value = 0.01929905749994466
$ws.Range("H41").Value = 28.68829640061239
$ws.Range("G42").Value = 0.1170014508066673
$ws.Range("H42").Value = 15.91747689128875
$ws.Range("G43").Value = 0.1238901287118265
$ws.Range("H43").Value = 3.117600384039898
$ws.Range("G44").Value = 0.01753072997119498
$ws.Range("H44").Value = -50.87723377176749
$ws.Range("G45").Value = 0.02799954508216979
$ws.Range("H45").Value = 71.03999507054091
$ws.Range("G46").Value = 0.06119905471358636
$ws.Range("H46").Value = 68.87846457427807
$ws.Range("G47").Value = 0.06949628160693959
$ws.Range("H47").Value = 37.77924100754851
$ws.Range("G48").Value = 0.02206608365955405
$ws.Range("H48").Value = -48.42119685002495
$ws.Range("G49").Value = 0.03853070655171578
$ws.Range("H49").Value = -44.54106350298274
$ws.Range("G50").Value = 0.03112198766653721
$ws.Range("H50").Value = 80.18014541514376
$ws.Range("G51").Value = 0.04304465516471798
$ws.Range("H51").Value = 121.0852285483004
$ws.Range("G52").Value = -0.1076335459548217
$ws.Range("H52").Value = -3.973215270702664
$ws.Range("G53").Value = -0.08182742571573319
$ws.Range("H53").Value = 11.39903805319407
$ws.Range("G54").Value = 0.09172326355517363
$ws.Range("H54").Value = 25.43753503113061
$ws.Range("G55").Value = 0.07790929710140275
$ws.Range("H55").Value = 25.75821286146874
$ws.Range("G56").Value = 0.03163388128049504
$ws.Range("H56").Value = -9.589622577774325
$ws.Range("G57").Value = -0.00755322800037414
$ws.Range("H57").Value = -230.8254058975033
$ws.Range("G58").Value = 0.03208914107160291
$ws.Range("H58").Value = 28.3021285185804
$ws.Range("G59").Value = 0.03752646890781916
$ws.Range("H59").Value = 58.48246116872795
$ws.Range("G60").Value = 0.03294721855803148
$ws.Range("H60").Value = 1.555238046642488
$ws.Range("G61").Value = 0.006224637330524719
$ws.Range("H61").Value = -50.82445339462709
$ws.Range("G62").Value = 0.06769913375723428
$ws.Range("H62").Value = 12.1540716770478
$ws.Range("G63").Value = 0.05590616199805591
$ws.Range("H63").Value = 71.54658222721679
$ws.Range("G64").Value = 0.01848840694080496
$ws.Range("H64").Value = -54.37910387803518
$ws.Range("G65").Value = 0.04976151977798025
$ws.Range("H65").Value = -11.23829016523639
$ws.Range("G66").Value = 0.09360367629722982
$ws.Range("H66").Value = 0.05298720819791147
$ws.Range("G67").Value = 0.0988613738629707
$ws.Range("H67").Value = -14.36618701476395
$ws.Range("G68").Value = -0.0167433095583018
$ws.Range("H68").Value = 51.95673580421757
$ws.Range("G69").Value = -0.00238756824422157
$ws.Range("H69").Value = 88.74948262397045
$ws.Range("G70").Value = 0.09136055483698134
$ws.Range("H70").Value = -1.377109944218387
$ws.Range("G71").Value = 0.1017772330056154
$ws.Range("H71").Value = 11.58708056324186
$ws.Range("G72").Value = -0.04829176454060918
$ws.Range("H72").Value = 13.88929935230405
$ws.Range("G73").Value = -0.06905767654755533
$ws.Range("H73").Value = 6.378454625037357
$ws.Range("G74").Value = 0.1131435076837879
$ws.Range("H74").Value = 13.20372437285388
$ws.Range("G75").Value = 0.1259883819826617
$ws.Range("H75").Value = 29.34617945286561
$ws.Range("G76").Value = -0.007708053796423714
$ws.Range("H76").Value = -130.1442058062402
$ws.Range("G77").Value = 0.01109402432663172
$ws.Range("H77").Value = -21.36939005002051
$ws.Range("G78").Value = 0.111303913918791
$ws.Range("H78").Value = 73.16297939009576
$ws.Range("G79").Value = 0.1065079716764543
$ws.Range("H79").Value = 38.83791190068428
$ws.Range("G80").Value = -0.1847227521994877
$ws.Range("H80").Value = -11.54189618286048
$ws.Range("G81").Value = -0.1432004410388426
$ws.Range("H81").Value = 31.83917943239663
$ws.Range("G82").Value = 0.1257243782375966
$ws.Range("H82").Value = 9.615954205491613
$ws.Range("G83").Value = 0.1865949915671559
$ws.Range("H83").Value = 4.839123493920554
$ws.Range("G84").Value = 0.06514484374669685
$ws.Range("H84").Value = 173.285522519327
$ws.Range("G85").Value = 0.06839527750644968
$ws.Range("H85").Value = 11.07467831954418
